# Loan RBI, Variable Instalments
# Insert a new (currently blank) "Variable Instalments" column into the
# "Repayment schedule" sheet, right before the existing "Late" column,
# and leave the "Repayment schedule" tab as the active/selected sheet
# (with cell L15 selected on it) instead of "Transactions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at N (shifts old N..P -> O..Q).  Excel's default
# insert behaviour copies formatting from the column to the left (M),
# which is why the new column picks up the same custom width as M.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and select cell L15 on it.
$ws.Activate()
$ws.Range("L15").Select() | Out-Null
